$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2-498).
# All of these values shift from 45181 (2023-09-12) to 45182 (2023-09-13).
$ws.Range("C2:C498").Value = 45182
